# Slide 8 (sldId 261) - shape id 12 "TextBox 11": append a sentence to the
# existing paragraph and let the autofit textbox grow to fit the new text.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(11)

$sh.TextFrame.TextRange.Text = "Les deux schémas étant d’ordre 2, il est attendu que l’ordre de convergence soit également d’ordre 2. De manière similaire à précédemment, on peut vérifier ce code.  "

# The shape uses <a:spAutoFit/>, so PowerPoint grows the box to fit the
# extra line of text; reproduce the resulting height (EMU 369332 -> 646331).
$sh.Height = 50.89225
